# Apply updated market-price snapshot values (currentAveragePrice* / Leve*Profit*
# columns H:N) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR leve-profit tables.
# Source data refreshed by the scheduled price-scrape runner; only cached
# values change, no formulas/layout are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 360.2857
$ws.Range("I19").Value = 364.63635
$ws.Range("K19").Value = 364.63635
$ws.Range("M19").Value = -189.63635
# Row 74
$ws.Range("H74").Value = 4505
$ws.Range("I74").Value = 4450
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 4450
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -3514
$ws.Range("N74").Value = -6872
# Row 77
$ws.Range("H77").Value = 4505
$ws.Range("I77").Value = 4450
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 22250
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -17570
$ws.Range("N77").Value = -34360
# Row 80
$ws.Range("H80").Value = 623.3182
$ws.Range("I80").Value = 544.5714
$ws.Range("J80").Value = 761.125
$ws.Range("K80").Value = 1633.7142
$ws.Range("L80").Value = 2283.375
$ws.Range("M80").Value = -635.7142000000001
$ws.Range("N80").Value = -4279.375
# Row 83
$ws.Range("H83").Value = 623.3182
$ws.Range("I83").Value = 544.5714
$ws.Range("J83").Value = 761.125
$ws.Range("K83").Value = 4901.1426
$ws.Range("L83").Value = 6850.125
$ws.Range("M83").Value = 90.85739999999987
$ws.Range("N83").Value = -16834.125
# Row 92
$ws.Range("H92").Value = 15556215
$ws.Range("I92").Value = 2778412.5
$ws.Range("K92").Value = 2778412.5
$ws.Range("M92").Value = -2777164.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5351.772
$ws.Range("I32").Value = 3955.7778
$ws.Range("K32").Value = 3955.7778
$ws.Range("M32").Value = -3668.7778
# Row 97
$ws.Range("H97").Value = 1323
$ws.Range("I97").Value = 1003.3333
$ws.Range("J97").Value = 1706.6
$ws.Range("K97").Value = 1003.3333
$ws.Range("L97").Value = 1706.6
$ws.Range("M97").Value = -507.3333
$ws.Range("N97").Value = -2698.6
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
# Row 122
$ws.Range("H122").Value = 1351701.5
$ws.Range("I122").Value = 2565847
$ws.Range("J122").Value = 2651
$ws.Range("K122").Value = 7697541
$ws.Range("L122").Value = 7953
$ws.Range("M122").Value = -7695091
$ws.Range("N122").Value = -12853
# Row 123
$ws.Range("H123").Value = 44427.668
$ws.Range("J123").Value = 44427.668
$ws.Range("L123").Value = 44427.668
$ws.Range("N123").Value = -54227.668

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 98
$ws.Range("H98").Value = 52780
$ws.Range("J98").Value = 52780
$ws.Range("L98").Value = 52780
$ws.Range("N98").Value = -57272
# Row 99
$ws.Range("H99").Value = 6483.3335
$ws.Range("I99").Value = 1800
$ws.Range("J99").Value = 7420
$ws.Range("K99").Value = 1800
$ws.Range("L99").Value = 7420
$ws.Range("M99").Value = -302
$ws.Range("N99").Value = -10416
# Row 126
$ws.Range("H126").Value = 6483.3335
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 7420
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 22260
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -27200

$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 12000
$ws.Range("J39").Value = 12000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -13064
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
# Row 122
$ws.Range("H122").Value = 3828852
$ws.Range("I122").Value = 3602422
$ws.Range("J122").Value = 4168497
$ws.Range("K122").Value = 10807266
$ws.Range("L122").Value = 12505491
$ws.Range("M122").Value = -10804816
$ws.Range("N122").Value = -12510391
# Row 126
$ws.Range("H126").Value = 7142.684
$ws.Range("I126").Value = 9992.583000000001
$ws.Range("J126").Value = 2257.1428
$ws.Range("K126").Value = 29977.749
$ws.Range("L126").Value = 6771.428400000001
$ws.Range("M126").Value = -27507.749
$ws.Range("N126").Value = -11711.4284
# Row 132
$ws.Range("H132").Value = 3447.383
$ws.Range("I132").Value = 3198.1667
$ws.Range("J132").Value = 3707.4348
$ws.Range("K132").Value = 9594.500100000001
$ws.Range("L132").Value = 11122.3044
$ws.Range("M132").Value = -7064.500100000001
$ws.Range("N132").Value = -16182.3044

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 125051750
$ws.Range("I93").Value = 200000
$ws.Range("J93").Value = 166669000
$ws.Range("K93").Value = 200000
$ws.Range("L93").Value = 166669000
$ws.Range("M93").Value = -198752
$ws.Range("N93").Value = -166671496
# Row 122
$ws.Range("H122").Value = 3970453.2
$ws.Range("I122").Value = 5104282.5
$ws.Range("J122").Value = 2050
$ws.Range("K122").Value = 15312847.5
$ws.Range("L122").Value = 6150
$ws.Range("M122").Value = -15310397.5
$ws.Range("N122").Value = -11050
# Row 132
$ws.Range("H132").Value = 10574098
$ws.Range("I132").Value = 12385967
$ws.Range("J132").Value = 4866.3335
$ws.Range("K132").Value = 37157901
$ws.Range("L132").Value = 14599.0005
$ws.Range("M132").Value = -37155371
$ws.Range("N132").Value = -19659.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 56
$ws.Range("H56").Value = 35907
$ws.Range("I56").Value = 5000
$ws.Range("J56").Value = 46209.332
$ws.Range("K56").Value = 5000
$ws.Range("L56").Value = 46209.332
$ws.Range("M56").Value = -4286
$ws.Range("N56").Value = -47637.332
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
# Row 100
$ws.Range("H100").Value = 12921.125
$ws.Range("I100").Value = 25393.5
$ws.Range("J100").Value = 448.75
$ws.Range("K100").Value = 50787
$ws.Range("L100").Value = 897.5
$ws.Range("M100").Value = -50246
$ws.Range("N100").Value = -1979.5
# Row 132
$ws.Range("H132").Value = 1144.9667
$ws.Range("I132").Value = 772.7917
$ws.Range("J132").Value = 2633.6667
$ws.Range("K132").Value = 2318.3751
$ws.Range("L132").Value = 7901.000100000001
$ws.Range("M132").Value = 211.6248999999998
$ws.Range("N132").Value = -12961.0001
